$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# (46074 -> 46075) for every data row (rows 2 through 408).
$ws.Range("C2:C408").Value = 46075
